$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 111528203
$ws.Range("AC10").Value = ""
$ws.Range("Q10").Value = 467745.6122397452
$ws.Range("R10").Value = 6875429.258361855
$ws.Range("A11").Value = 111527806
$ws.Range('P11').Value = 'Fläcksberget V, Hjd'
$ws.Range("Q11").Value = 467615.2905344999
$ws.Range("R11").Value = 6875426.740629551
$ws.Range("A12").Value = 111528365
$ws.Range('P12').Value = 'Fläcksberget, Hjd'
$ws.Range("Q12").Value = 467795.2212022893
$ws.Range("R12").Value = 6875452.272210476
$ws.Range("A13").Value = 111528980
$ws.Range("Q13").Value = 467799.8074815667
$ws.Range("R13").Value = 6875539.119922069
$ws.Range("A14").Value = 111909536
$ws.Range("B14").Value = 77267
$ws.Range("E14").Value = 6446
$ws.Range('F14').Value = 'Kolflarnlav'
$ws.Range('G14').Value = 'Carbonicola anthracophila'
$ws.Range('H14').Value = '(Nyl.) Bendiksby & Timdal'
$ws.Range('P14').Value = 'Fläcksberget, Hjd'
$ws.Range("Q14").Value = 467891.3929605001
$ws.Range("R14").Value = 6875425.059267788
$ws.Range("A15").Value = 111908768
$ws.Range("B15").Value = 96348
$ws.Range('D15').Value = 'VU'
$ws.Range("E15").Value = 220787
$ws.Range('F15').Value = 'Knärot'
$ws.Range('G15').Value = 'Goodyera repens'
$ws.Range('H15').Value = '(L.) R. Br.'
$ws.Range('I15').Value = "'1"
$ws.Range("Q15").Value = 467911.8445363804
$ws.Range("R15").Value = 6875299.456096188
$ws.Range("A16").Value = 111908364
$ws.Range("AC16").Value = ""
$ws.Range("B16").Value = 90660
$ws.Range('D16').Value = 'NT'
$ws.Range("E16").Value = 4362
$ws.Range('F16').Value = 'Blå taggsvamp'
$ws.Range('G16').Value = 'Hydnellum caeruleum'
$ws.Range('H16').Value = '(Hornem.) P.Karst.'
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range('P16').Value = 'Gröbäcken, Hjd'
$ws.Range("Q16").Value = 467724.2196293612
$ws.Range("R16").Value = 6874811.291555981
$ws.Range("A17").Value = 111909174
$ws.Range("B17").Value = 77267
$ws.Range('D17').Value = 'NT'
$ws.Range("E17").Value = 6446
$ws.Range('F17').Value = 'Kolflarnlav'
$ws.Range('G17').Value = 'Carbonicola anthracophila'
$ws.Range('H17').Value = '(Nyl.) Bendiksby & Timdal'
$ws.Range("I17").Value = ""
$ws.Range("Q17").Value = 467989.0228066717
$ws.Range("R17").Value = 6875352.744105402
$ws.Range("A19").Value = 111908700
$ws.Range('AC19').Value = 'Förekomst av doftticka i avverkningsanmält område.'
$ws.Range("B19").Value = 89965
$ws.Range('D19').Value = 'VU'
$ws.Range("E19").Value = 760
$ws.Range('F19').Value = 'Doftticka'
$ws.Range('G19').Value = 'Haploporus odorus'
$ws.Range('H19').Value = '(Sommerf.) Bondartsev & Singer'
$ws.Range('I19').Value = "'6"
$ws.Range('J19').Value = 'fruktkroppar'
$ws.Range("Q19").Value = 467921.7931363151
$ws.Range("R19").Value = 6875306.87748003
$ws.Range("A21").Value = 112014208
$ws.Range("Q21").Value = 467418.043506761
$ws.Range("R21").Value = 6875312.610613029
$ws.Range("A22").Value = 112015011
$ws.Range("Q22").Value = 467389.9660160011
$ws.Range("R22").Value = 6875327.91063729
$ws.Range("A23").Value = 112014229
$ws.Range("B23").Value = 90682
$ws.Range("E23").Value = 2059
$ws.Range('F23').Value = 'Skrovlig taggsvamp'
$ws.Range('G23').Value = 'Hydnellum scabrosum'
$ws.Range('H23').Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q23").Value = 467427.230114766
$ws.Range("R23").Value = 6875289.506732536
$ws.Range("A24").Value = 112014423
$ws.Range("B24").Value = 90658
$ws.Range("E24").Value = 4361
$ws.Range('F24').Value = 'Orange taggsvamp'
$ws.Range('G24').Value = 'Hydnellum aurantiacum'
$ws.Range('H24').Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("Q24").Value = 467430.0274016621
$ws.Range("R24").Value = 6875237.811246304
$ws.Range("A25").Value = 112014177
$ws.Range("B25").Value = 90689
$ws.Range('D25').Value = 'NT'
$ws.Range("E25").Value = 5966
$ws.Range('F25').Value = 'Motaggsvamp'
$ws.Range('G25').Value = 'Sarcodon squamosus'
$ws.Range('H25').Value = '(Schaeff.) Quél.'
$ws.Range("Q25").Value = 467389.9660160011
$ws.Range("R25").Value = 6875327.91063729
$ws.Range("A27").Value = 112014347
$ws.Range("B27").Value = 90678
$ws.Range('D27').Value = 'LC'
$ws.Range("E27").Value = 4366
$ws.Range('F27').Value = 'Skarp dropptaggsvamp'
$ws.Range('G27').Value = 'Hydnellum peckii'
$ws.Range('H27').Value = 'Banker'
$ws.Range("Q27").Value = 467430.0274016621
$ws.Range("R27").Value = 6875237.811246304
$ws.Range("A28").Value = 112014142
$ws.Range("B28").Value = 90666
$ws.Range('D28').Value = 'LC'
$ws.Range("E28").Value = 4364
$ws.Range('F28').Value = 'Dropptaggsvamp'
$ws.Range('G28').Value = 'Hydnellum ferrugineum'
$ws.Range('H28').Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q28").Value = 467442.7363991642
$ws.Range("R28").Value = 6875336.798642672
$ws.Range("A29").Value = 112014300
$ws.Range("B29").Value = 90689
$ws.Range('D29').Value = 'NT'
$ws.Range("E29").Value = 5966
$ws.Range('F29').Value = 'Motaggsvamp'
$ws.Range('G29').Value = 'Sarcodon squamosus'
$ws.Range('H29').Value = '(Schaeff.) Quél.'
$ws.Range("Q29").Value = 467415.4484496959
$ws.Range("R29").Value = 6875287.271149865
